# Daily "roll forward" update for the 剩余 (days-remaining) tracker sheet.
#
# Each data row tracks a cycle: D = 总天 (total days in cycle), F = 开始时间
# (cycle start date, stored as a plain yyyymmdd number, e.g. 20260209), and
# E = 剩余 (days remaining) = D - (today - F).
#
# This script advances "today" by one day (the workbook was last computed as
# of 2026-02-12; this run recomputes it as of 2026-02-13) and rewrites E (and,
# when a cycle has just run out, F) for every data row accordingly:
#   - if the recomputed remaining days is still > 0, only E changes (drops by
#     one day elapsed);
#   - if the cycle would hit zero/negative, the cycle "resets": F becomes the
#     new "today" and E goes back to the full D (a fresh cycle just started).
#
# Row 35 (sheet row 36) has a corrupted 9-digit start date (202510929) that
# cannot be parsed as yyyymmdd, so — matching the source data — it is left
# untouched, exactly as the upstream automation that produced this data
# apparently skipped it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "as of" date for this run (one day after the previous run).
$newToday = 20260213

# Find the last used row in column A (行号) so this keeps working if rows are
# appended/removed.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

# Use the Application's own formula engine (via Evaluate) to do yyyymmdd ->
# date-serial conversion, so month/year rollovers are handled correctly
# without reimplementing a calendar and without touching any worksheet cell.
$newTodaySerial = $excel.Evaluate("DATE(LEFT($newToday,4),MID($newToday,5,2),RIGHT($newToday,2))")

for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $fVal = $fCell.Value2
    $total = $dCell.Value2

    if ($null -eq $fVal -or $fVal -eq "" -or $null -eq $total -or $total -eq "") {
        continue
    }

    $fStr = [string][int64]$fVal
    if ($fStr.Length -ne 8) {
        # Malformed start date (e.g. "202510929") - skip, same as source data.
        continue
    }

    $fDateSerial = $excel.Evaluate("DATE(LEFT($fVal,4),MID($fVal,5,2),RIGHT($fVal,2))")

    $newRemaining = $total - ($newTodaySerial - $fDateSerial)

    if ($newRemaining -le 0) {
        # Cycle elapsed: start a fresh one as of today.
        $eCell.Value2 = $total
        $fCell.Value2 = $newToday
    } else {
        $eCell.Value2 = $newRemaining
    }
}
